# mute_contact.xlsx uplift: drop the redundant "form_id" setting (the
# form_id is now derived elsewhere by the newer pyxform/cht-conf tooling)
# and clean up a stray leftover "NO_LABEL" value on the survey sheet.

$wb = $excel.ActiveWorkbook
$wsSurvey   = $wb.Worksheets.Item("survey")
$wsSettings = $wb.Worksheets.Item("settings")

# --- settings sheet -------------------------------------------------
# Column B ("form_id") is being removed entirely, so column C onward
# (version / style / namespaces) shifts left by one. Before deleting the
# column, carry each header-comment's text one slot to the left so the
# remaining comments keep describing the right column after the shift,
# then drop the now-orphaned trailing comment.
$versionComment    = $wsSettings.Range("C1").Comment.Text()
$styleComment      = $wsSettings.Range("D1").Comment.Text()
$namespacesComment = $wsSettings.Range("E1").Comment.Text()

$wsSettings.Range("B1").Comment.Text($versionComment)
$wsSettings.Range("C1").Comment.Text($styleComment)
$wsSettings.Range("D1").Comment.Text($namespacesComment)
$wsSettings.Range("E1").Comment.Delete()

$wsSettings.Columns("B").Delete()
$wsSettings.Range("B1").Select()

# --- survey sheet -----------------------------------------------------
# Clear the stray "NO_LABEL" appearance value that was left in C3.
$wsSurvey.Range("C3").Clear()
$wsSurvey.Range("C3").Select()
